$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "61.612.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -2.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.007.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -2.79%  "
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "539.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "136.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.25%  "
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.001.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -2.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -4.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000223"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "34.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.479.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -2.90%  "
$ws.Range("E16").Value2 = "  -0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "61.568.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.991.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "470.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -3.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.680"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -3.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "80.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "12.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "7.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -5.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.03%  "
$ws.Range("B30").Value2 = "ImmutableX"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -0.38%  "
$ws.Range("B31").Value2 = "EthereumClassic"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "25.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -1.92%  "
$ws.Range("B32").Value2 = "Mantle"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +3.04%  "
$ws.Range("B33").Value2 = "NEARProtocol"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -4.32%  "
$ws.Range("B35").Value2 = "OKB"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "55.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "459.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -6.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "3.173.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -1.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.0388"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "8.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -6.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "27.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +11.93%  "
$ws.Range("E45").Value2 = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.247"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -3.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "119.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.109"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0₃0502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -8.07%  "
$ws.Range("B51").Value2 = "ThetaToken"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -1.11%  "
